$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61 (pushes old rows 61-190 down to 62-191, and
# extends the used range to R191, mirroring a weekly-refresh prepend of
# one new daily price record).
$ws.Rows("61").Insert()

# Populate the newly inserted row 61 with the new record.
$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44536
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = 100112012
$ws.Range("G61").Value = "Espinaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 400
$ws.Range("L61").Value = 500
$ws.Range("M61").Value = 450
$ws.Range("N61").Value = "$/atado 300 a 500 gramos"
$ws.Range("O61").Value = "Provincia del Elquí"
$ws.Range("P61").Value = 900
$ws.Range("Q61").Value = 0.5
$ws.Range("R61").Value = "Hortaliza"
